$wb = $excel.ActiveWorkbook

# Rename the sheets
$wb.Worksheets.Item(1).Name = "Thermal Battery"
$wb.Worksheets.Item(2).Name = "Greenhouse"

# Add a new bold "Item" header cell at A18 on the Greenhouse sheet
$ws2 = $wb.Worksheets.Item("Greenhouse")
$ws2.Range("A18").Value = "Item"
$ws2.Range("A18").Font.Bold = $true
